{"js": "// Update the division-problem worksheet numbers in the single table.\n// Each \"data\" row of the table holds five \"a\u00f7b=\" expressions; the rows\n// are matched by their current (old) values so the edit is resilient to\n// exactly which row indices hold data.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Map old row (as it appears before the edit) -> new row.\nconst rowReplacements = [\n  { oldRow: [\"53\u00f75=\", \"27\u00f75=\", \"18\u00f77=\", \"43\u00f76=\", \"79\u00f78=\"], newRow: [\"25\u00f73=\", \"35\u00f79=\", \"55\u00f77=\", \"34\u00f74=\", \"41\u00f72=\"] },\n  { oldRow: [\"83\u00f78=\", \"88\u00f74=\", \"46\u00f79=\", \"78\u00f77=\", \"55\u00f79=\"], newRow: [\"48\u00f74=\", \"18\u00f78=\", \"18\u00f78=\", \"28\u00f79=\", \"28\u00f79=\"] },\n  { oldRow: [\"86\u00f75=\", \"89\u00f75=\", \"97\u00f73=\", \"84\u00f72=\", \"78\u00f77=\"], newRow: [\"89\u00f77=\", \"67\u00f76=\", \"19\u00f75=\", \"66\u00f73=\", \"89\u00f72=\"] },\n  { oldRow: [\"41\u00f79=\", \"80\u00f77=\", \"96\u00f75=\", \"10\u00f79=\", \"55\u00f75=\"], newRow: [\"52\u00f73=\", \"52\u00f74=\", \"43\u00f74=\", \"65\u00f78=\", \"77\u00f76=\"] },\n  { oldRow: [\"13\u00f77=\", \"33\u00f72=\", \"62\u00f78=\", \"83\u00f75=\", \"38\u00f74=\"], newRow: [\"54\u00f74=\", \"94\u00f73=\", \"14\u00f75=\", \"15\u00f77=\", \"60\u00f78=\"] },\n];\n\nconst values = table.values;\nfor (const row of values) {\n  const match = rowReplacements.find((r) => r.oldRow.every((v, i) => v === row[i]));\n  if (match) {\n    for (let i = 0; i < match.newRow.length; i++) {\n      row[i] = match.newRow[i];\n    }\n  }\n}\n\ntable.values = values;\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet numbers in the single table.\n# Each \"data\" row of the table holds five \"a\u00f7b=\" expressions; rows are\n# matched by their current (old) values so the edit is resilient to\n# exactly which row indices hold data.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowReplacements = @(\n    @{ Old = @(\"53\u00f75=\", \"27\u00f75=\", \"18\u00f77=\", \"43\u00f76=\", \"79\u00f78=\"); New = @(\"25\u00f73=\", \"35\u00f79=\", \"55\u00f77=\", \"34\u00f74=\", \"41\u00f72=\") },\n    @{ Old = @(\"83\u00f78=\", \"88\u00f74=\", \"46\u00f79=\", \"78\u00f77=\", \"55\u00f79=\"); New = @(\"48\u00f74=\", \"18\u00f78=\", \"18\u00f78=\", \"28\u00f79=\", \"28\u00f79=\") },\n    @{ Old = @(\"86\u00f75=\", \"89\u00f75=\", \"97\u00f73=\", \"84\u00f72=\", \"78\u00f77=\"); New = @(\"89\u00f77=\", \"67\u00f76=\", \"19\u00f75=\", \"66\u00f73=\", \"89\u00f72=\") },\n    @{ Old = @(\"41\u00f79=\", \"80\u00f77=\", \"96\u00f75=\", \"10\u00f79=\", \"55\u00f75=\"); New = @(\"52\u00f73=\", \"52\u00f74=\", \"43\u00f74=\", \"65\u00f78=\", \"77\u00f76=\") },\n    @{ Old = @(\"13\u00f77=\", \"33\u00f72=\", \"62\u00f78=\", \"83\u00f75=\", \"38\u00f74=\"); New = @(\"54\u00f74=\", \"94\u00f73=\", \"14\u00f75=\", \"15\u00f77=\", \"60\u00f78=\") }\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $cols = $t.Columns.Count\n    $current = @()\n    for ($c = 1; $c -le $cols; $c++) {\n        $txt = $t.Cell($r, $c).Range.Text\n        $current += ($txt -replace \"[\\r\\x07]\", \"\")\n    }\n\n    foreach ($mapping in $rowReplacements) {\n        $isMatch = $true\n        for ($c = 0; $c -lt $cols; $c++) {\n            if ($current[$c] -ne $mapping.Old[$c]) {\n                $isMatch = $false\n                break\n            }\n        }\n        if ($isMatch) {\n            for ($c = 1; $c -le $cols; $c++) {\n                $t.Cell($r, $c).Range.Text = $mapping.New[$c - 1]\n            }\n            break\n        }\n    }\n}\n"}
